$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the a1a03878... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-19 06:44:18"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime
# for the a1a03878... row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-19 06:44:13"
$wsZhCn.Range("K4").Value = "2016-08-19 06:44:31"

# de-de sheet: update Correspond Handoff Datetime and Correspond Handback DateTime
# for the a1a03878... row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-19 06:44:18"
$wsDeDe.Range("K4").Value = "2016-08-19 06:44:40"
